# The "Pens" transaction (row 14) on the Transactions sheet is being
# removed. Deleting the whole row shifts rows 15-19 up by one, which
# matches the target: row 15 (Powerwheel for Lucas) becomes row 14,
# row 16 (New computer) becomes row 15, row 17 (Tire for dirtbike)
# becomes row 16, row 18 (New bib for tire) becomes row 17, and
# row 19 (Notebooks for work) becomes row 18.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Rows.Item(14).Delete()
